# ontoDog_input.xlsx update:
#  - add a new "age measurement datum" row (OBI_0001167) to Sheet1
#  - update the sheet view's scroll position / selected cell to match
#    the author's on-disk view state at commit time

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# New row 53: IRI | label | "y" (include in view) | (blank user label) | "yes" (include all children)
$ws.Range("A53").Value = "http://purl.obolibrary.org/obo/OBI_0001167"
$ws.Range("B53").Value = "age measurement datum"
$ws.Range("C53").Value = "y"
$ws.Range("E53").Value = "yes"

# Restore the view: scrolled to row 37, active/selected cell B57
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B57").Select()
